$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginTestData")
$ws.Activate()

# Resize/reposition the workbook window (mirrors the bookViews/workbookView
# xWindow/yWindow/windowWidth/windowHeight change in the saved file).
$aw = $excel.ActiveWindow
$aw.WindowState = "xlNormal"
$aw.Width  = 15570
$aw.Height = 4830

# Scroll the sheet so column D is the left-most visible column, then move the
# selection from J2 to E2 (new topLeftCell="D1" + selection activeCell/sqref).
$aw.ScrollColumn = 4
$aw.ScrollRow = 1
$ws.Range("E2").Select()

# New data point entered in K2.
$ws.Range("K2").Value = 1
